$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update the "Steps" (G3) and "Expected Behaviour" (H3) cells for the
# VT200_0681 test case: the trailing screenshot step/validation was
# replaced with a plain text-contains check, and the now-unused
# "validate5" screenshot block was dropped.
$ws.Range("G3").Value = "wait(5);`nvalidate1;`nlink_Click(notification_test_link);`nvalidate2;`nSelectTestToRun(VT200_0681_string);`nClickRunTest(runtest_top_xpath);`nwait(2);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(3);`nCheckUITextContains(This_is_a_pop_up_for_hide);`nwait(10);`nCheckUITextContains(hidepopup);"

$ws.Range("H3").Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Notification JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0681`n};"

# The shorter text now wraps into fewer lines, so the row shrinks.
$ws.Rows.Item(3).RowHeight = 166.5

# Move the active selection.
$ws.Range("D2").Select()
